$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff / Handback datetimes for first data row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 20:53:09"
$wsZh.Range("H2").Value = "2016-03-19 20:53:34"

# de-de sheet: update Correspond Handoff / Handback datetimes for first data row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 20:53:12"
$wsDe.Range("H2").Value = "2016-03-19 20:53:39"
